{"js": "// Load all paragraphs in the document body so we can locate the two\n// anchor paragraphs by their text/style (more robust than raw index,\n// in case surrounding content shifts).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text,style\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// 1) Find the (first) paragraph that sits right after the \"\u0422\u0435\u043c\u0430 ...\"\n//    title paragraph and currently holds just a single space \u2014 this is\n//    the paragraph the new \"\u0412\u0430\u0440\u0456\u0430\u043d\u0442 \u211619\" line must be inserted before.\nlet targetIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.trim() === \"\" && i > 0 && items[i - 1].text.indexOf(\"\u0422\u0435\u043c\u0430\") !== -1) {\n    targetIndex = i;\n    break;\n  }\n}\nif (targetIndex === -1) {\n  throw new Error(\"Could not locate the paragraph following the '\u0422\u0435\u043c\u0430' line.\");\n}\n\n// Insert a new paragraph with the variant number directly before it;\n// it inherits that paragraph's formatting (centered, no first-line\n// indent, contextual spacing) automatically.\nitems[targetIndex].insertParagraph(\"\u0412\u0430\u0440\u0456\u0430\u043d\u0442 \u211619\", \"Before\");\n\n// 2) Find and remove the empty paragraph styled \"P\" that sits between\n//    the \"2024\" paragraph and the \"\u0422\u0430\u043a \u0431\u043e \u0411\u041e\u0413...\" Quote paragraph.\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].style === \"P\" && items[i].text === \"\" &&\n      i > 0 && items[i - 1].text.indexOf(\"2024\") !== -1) {\n    items[i].delete();\n    break;\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Locate the paragraph that sits right after the \"\u0422\u0435\u043c\u0430 ...\" title\n#    paragraph and currently contains just a single space \u2014 insert the\n#    new \"\u0412\u0430\u0440\u0456\u0430\u043d\u0442 \u211619\" paragraph immediately before it.\n$targetIndex = -1\nfor ($i = 2; $i -le $d.Paragraphs.Count; $i++) {\n    $cur = $d.Paragraphs.Item($i)\n    $prev = $d.Paragraphs.Item($i - 1)\n    if ($cur.Range.Text.Trim().Length -eq 0 -and $prev.Range.Text.Contains(\"\u0422\u0435\u043c\u0430\")) {\n        $targetIndex = $i\n        break\n    }\n}\nif ($targetIndex -eq -1) {\n    throw \"Could not locate the paragraph following the '\u0422\u0435\u043c\u0430' line.\"\n}\n\n$targetRange = $d.Paragraphs.Item($targetIndex).Range\n$targetRange.InsertParagraphBefore()\n$newPara = $d.Paragraphs.Item($targetIndex)\n$newPara.Range.Text = \"\u0412\u0430\u0440\u0456\u0430\u043d\u0442 \u211619\"\n\n# 2) Locate and remove the empty paragraph styled \"P\" that sits between\n#    the \"2024\" paragraph and the \"\u0422\u0430\u043a \u0431\u043e \u0411\u041e\u0413...\" Quote paragraph.\n$deleteIndex = -1\nfor ($i = 2; $i -le $d.Paragraphs.Count; $i++) {\n    $cur = $d.Paragraphs.Item($i)\n    $prev = $d.Paragraphs.Item($i - 1)\n    if ($cur.Style.NameLocal -eq \"P\" -and $cur.Range.Text.Trim().Length -eq 0 -and $prev.Range.Text.Contains(\"2024\")) {\n        $deleteIndex = $i\n        break\n    }\n}\nif ($deleteIndex -eq -1) {\n    throw \"Could not locate the empty 'P' paragraph following the '2024' line.\"\n}\n\n$d.Paragraphs.Item($deleteIndex).Range.Delete()\n"}
